# download articles with pandoc title blocks
#
# Rebuilds the opening "title block" of the article:
#   - title paragraph:  Heading1 -> Title style, text split into one run
#     per word/space token (as a pandoc-style docx writer would emit)
#   - byline paragraph: bold "By Dorothy Day" -> Authors-styled paragraph
#     with just "Dorothy Day", again split into one run per token
#   - the stray bookmark wrapping the old title is dropped along with it
#
# Word's Range.Text/InsertBefore/InsertAfter happily coalesces adjacent
# same-formatted runs back into a single <w:r>, so to land the exact
# one-run-per-token shape pandoc produces we build the paragraph XML
# ourselves and drop it in with Range.InsertXML (the same mechanism
# Word uses under the hood for XML-fragment paste).

$d = $word.ActiveDocument

# The old title sat inside a "on-pilgrimage---november-1953" bookmark; the
# new title block drops it. Best-effort removal up front (harmless if the
# host doesn't surface bookmarks at all - we still rebuild the paragraphs
# below regardless).
try {
    $bmName = "on-pilgrimage---november-1953"
    if ($d.Bookmarks.Exists($bmName)) {
        $d.Bookmarks.Item($bmName).Delete()
    }
} catch {
}

function New-RunXml([string]$text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    return '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
}

function New-TokenRuns([string]$phrase) {
    $words = $phrase.Split(" ")
    $runs = ""
    for ($i = 0; $i -lt $words.Count; $i++) {
        if ($i -gt 0) {
            $runs += New-RunXml " "
        }
        $runs += New-RunXml $words[$i]
    }
    return $runs
}

$titleRuns = New-TokenRuns "On Pilgrimage - November 1953"
$authorRuns = New-TokenRuns "Dorothy Day"

$packageOpen = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$packageClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titlePara = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>'
$authorPara = '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>'

$xmlChunk = $packageOpen + $titlePara + $authorPara + $packageClose

# Paragraph 1 is the "On Pilgrimage - November 1953" Heading1 title,
# paragraph 2 is the bold "By Dorothy Day" byline directly under it.
$titleParagraph = $d.Paragraphs.Item(1)
$bylineParagraph = $d.Paragraphs.Item(2)
$span = $d.Range($titleParagraph.Range.Start, $bylineParagraph.Range.End)
$span.InsertXML($xmlChunk)

Write-Output "Title/byline block rebuilt."
